# Reapply the latest scheduled market-data refresh to the Leve profit sheets.
#
# Each worksheet ("ALC", "ARM", ... one per crafting job) has columns:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
# which are refreshed from the market-board snapshot. This script pushes the
# refreshed numbers for the rows the runner flagged as stale, row by row.
# Cells the new snapshot has no data for are cleared instead of zeroed, and
# cells that newly gained data are created the same way.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 10028.529  # H32: 10087.353 -> 10028.529
$ws.Cells.Item(32, 10).Value = 9060.888999999999  # J32: 9172 -> 9060.888999999999
$ws.Cells.Item(32, 12).Value = 9060.888999999999  # L32: 9172 -> 9060.888999999999
$ws.Cells.Item(32, 14).Value = -9712.888999999999  # N32: -9824 -> -9712.888999999999

$ws.Cells.Item(43, 8).Value = 4885  # H43: 0 -> 4885
$ws.Cells.Item(43, 10).Value = 4885  # J43: 0 -> 4885
$ws.Cells.Item(43, 12).Value = 4885  # L43: 0 -> 4885
$ws.Cells.Item(43, 14).Value = -5023  # N43: (blank) -> -5023

$ws.Cells.Item(80, 8).Value = 1954.7142  # H80: 1960 -> 1954.7142
$ws.Cells.Item(80, 9).Value = 1859  # I80: 2065 -> 1859
$ws.Cells.Item(80, 10).Value = 2007.8889  # J80: 1897 -> 2007.8889
$ws.Cells.Item(80, 11).Value = 5577  # K80: 6195 -> 5577
$ws.Cells.Item(80, 12).Value = 6023.6667  # L80: 5691 -> 6023.6667
$ws.Cells.Item(80, 13).Value = -4579  # M80: -5197 -> -4579
$ws.Cells.Item(80, 14).Value = -8019.6667  # N80: -7687 -> -8019.6667

$ws.Cells.Item(83, 8).Value = 1954.7142  # H83: 1960 -> 1954.7142
$ws.Cells.Item(83, 9).Value = 1859  # I83: 2065 -> 1859
$ws.Cells.Item(83, 10).Value = 2007.8889  # J83: 1897 -> 2007.8889
$ws.Cells.Item(83, 11).Value = 16731  # K83: 18585 -> 16731
$ws.Cells.Item(83, 12).Value = 18071.0001  # L83: 17073 -> 18071.0001
$ws.Cells.Item(83, 13).Value = -11739  # M83: -13593 -> -11739
$ws.Cells.Item(83, 14).Value = -28055.0001  # N83: -27057 -> -28055.0001

$ws.Cells.Item(137, 8).Value = 2599.25  # H137: 2673.5 -> 2599.25
$ws.Cells.Item(137, 9).Value = 2048.8333  # I137: 2026.2858 -> 2048.8333
$ws.Cells.Item(137, 10).Value = 3149.6667  # J137: 3579.6 -> 3149.6667
$ws.Cells.Item(137, 11).Value = 6146.499899999999  # K137: 6078.857400000001 -> 6146.499899999999
$ws.Cells.Item(137, 12).Value = 9449.000100000001  # L137: 10738.8 -> 9449.000100000001
$ws.Cells.Item(137, 13).Value = -3596.499899999999  # M137: -3528.857400000001 -> -3596.499899999999
$ws.Cells.Item(137, 14).Value = -14549.0001  # N137: -15838.8 -> -14549.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 343.8889  # H4: 378.33334 -> 343.8889
$ws.Cells.Item(4, 9).Value = 343.8889  # I4: 378.33334 -> 343.8889
$ws.Cells.Item(4, 11).Value = 343.8889  # K4: 378.33334 -> 343.8889
$ws.Cells.Item(4, 13).Value = -227.8889  # M4: -262.33334 -> -227.8889

$ws.Cells.Item(28, 8).Value = 12629.667  # H28: 9987.5 -> 12629.667
$ws.Cells.Item(28, 9).Value = 12629.667  # I28: 9987.5 -> 12629.667
$ws.Cells.Item(28, 11).Value = 12629.667  # K28: 9987.5 -> 12629.667
$ws.Cells.Item(28, 13).Value = -12437.667  # M28: -9795.5 -> -12437.667

$ws.Cells.Item(32, 8).Value = 3176.9  # H32: 2716.5588 -> 3176.9
$ws.Cells.Item(32, 9).Value = 3176.9  # I32: 2716.5588 -> 3176.9
$ws.Cells.Item(32, 11).Value = 3176.9  # K32: 2716.5588 -> 3176.9
$ws.Cells.Item(32, 13).Value = -2889.9  # M32: -2429.5588 -> -2889.9

$ws.Cells.Item(45, 8).Value = 3484.7144  # H45: 2967.125 -> 3484.7144
$ws.Cells.Item(45, 9).Value = 2398.6  # I45: 2333.8572 -> 2398.6
$ws.Cells.Item(45, 10).Value = 6200  # J45: 7400 -> 6200
$ws.Cells.Item(45, 11).Value = 2398.6  # K45: 2333.8572 -> 2398.6
$ws.Cells.Item(45, 12).Value = 6200  # L45: 7400 -> 6200
$ws.Cells.Item(45, 13).Value = -2021.6  # M45: -1956.8572 -> -2021.6
$ws.Cells.Item(45, 14).Value = -6954  # N45: -8154 -> -6954

$ws.Cells.Item(61, 8).Value = 1839.4  # H61: 1564.75 -> 1839.4
$ws.Cells.Item(61, 9).Value = 1724.25  # I61: 1497.909 -> 1724.25
$ws.Cells.Item(61, 11).Value = 1724.25  # K61: 1497.909 -> 1724.25
$ws.Cells.Item(61, 13).Value = -1512.25  # M61: -1285.909 -> -1512.25

$ws.Cells.Item(74, 8).Value = 3204.7727  # H74: 3373.05 -> 3204.7727
$ws.Cells.Item(74, 9).Value = 2808.111  # I74: 2968.875 -> 2808.111
$ws.Cells.Item(74, 11).Value = 2808.111  # K74: 2968.875 -> 2808.111
$ws.Cells.Item(74, 13).Value = -1934.111  # M74: -2094.875 -> -1934.111

$ws.Cells.Item(77, 8).Value = 3204.7727  # H77: 3373.05 -> 3204.7727
$ws.Cells.Item(77, 9).Value = 2808.111  # I77: 2968.875 -> 2808.111
$ws.Cells.Item(77, 11).Value = 14040.555  # K77: 14844.375 -> 14040.555
$ws.Cells.Item(77, 13).Value = -9672.555  # M77: -10476.375 -> -9672.555

$ws.Cells.Item(99, 8).Value = 12629.667  # H99: 9987.5 -> 12629.667
$ws.Cells.Item(99, 9).Value = 12629.667  # I99: 9987.5 -> 12629.667
$ws.Cells.Item(99, 11).Value = 12629.667  # K99: 9987.5 -> 12629.667
$ws.Cells.Item(99, 13).Value = -9634.666999999999  # M99: -6992.5 -> -9634.666999999999

$ws.Cells.Item(114, 8).Value = 30397.25  # H114: 30397.5 -> 30397.25
$ws.Cells.Item(114, 10).Value = 30397.25  # J114: 30397.5 -> 30397.25
$ws.Cells.Item(114, 12).Value = 30397.25  # L114: 30397.5 -> 30397.25
$ws.Cells.Item(114, 14).Value = -39075.25  # N114: -39075.5 -> -39075.25

$ws.Cells.Item(122, 8).Value = 0  # H122: 731.3333 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 1073 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 48 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 3219 -> 0
$ws.Cells.Item(122, 12).ClearContents()  # L122: 144 -> (blank)
$ws.Cells.Item(122, 13).ClearContents()  # M122: -769 -> (blank)
$ws.Cells.Item(122, 14).Value = 0  # N122: -5044 -> 0

$ws.Cells.Item(132, 8).Value = 6534.973  # H132: 6872.8237 -> 6534.973
$ws.Cells.Item(132, 9).Value = 6154.6772  # I132: 6333.6553 -> 6154.6772
$ws.Cells.Item(132, 10).Value = 8499.833000000001  # J132: 10000 -> 8499.833000000001
$ws.Cells.Item(132, 11).Value = 18464.0316  # K132: 19000.9659 -> 18464.0316
$ws.Cells.Item(132, 12).Value = 25499.499  # L132: 30000 -> 25499.499
$ws.Cells.Item(132, 13).Value = -15934.0316  # M132: -16470.9659 -> -15934.0316
$ws.Cells.Item(132, 14).Value = -30559.499  # N132: -35060 -> -30559.499

$ws.Cells.Item(136, 8).Value = 1839.4  # H136: 1564.75 -> 1839.4
$ws.Cells.Item(136, 9).Value = 1724.25  # I136: 1497.909 -> 1724.25
$ws.Cells.Item(136, 11).Value = 5172.75  # K136: 4493.727000000001 -> 5172.75
$ws.Cells.Item(136, 13).Value = -2622.75  # M136: -1943.727000000001 -> -2622.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 19999  # H26: 39999 -> 19999
$ws.Cells.Item(26, 9).Value = 19999  # I26: 39999 -> 19999
$ws.Cells.Item(26, 11).Value = 19999  # K26: 39999 -> 19999
$ws.Cells.Item(26, 13).Value = -19707  # M26: -39707 -> -19707

$ws.Cells.Item(86, 8).Value = 13001.333  # H86: 10776.267 -> 13001.333
$ws.Cells.Item(86, 9).Value = 18506.143  # I86: 20124.166 -> 18506.143
$ws.Cells.Item(86, 10).Value = 5294.6  # J86: 4544.3335 -> 5294.6
$ws.Cells.Item(86, 11).Value = 18506.143  # K86: 20124.166 -> 18506.143
$ws.Cells.Item(86, 12).Value = 5294.6  # L86: 4544.3335 -> 5294.6
$ws.Cells.Item(86, 13).Value = -17383.143  # M86: -19001.166 -> -17383.143
$ws.Cells.Item(86, 14).Value = -7540.6  # N86: -6790.3335 -> -7540.6

$ws.Cells.Item(89, 8).Value = 13001.333  # H89: 10776.267 -> 13001.333
$ws.Cells.Item(89, 9).Value = 18506.143  # I89: 20124.166 -> 18506.143
$ws.Cells.Item(89, 10).Value = 5294.6  # J89: 4544.3335 -> 5294.6
$ws.Cells.Item(89, 11).Value = 92530.715  # K89: 100620.83 -> 92530.715
$ws.Cells.Item(89, 12).Value = 26473  # L89: 22721.6675 -> 26473
$ws.Cells.Item(89, 13).Value = -86914.715  # M89: -95004.83 -> -86914.715
$ws.Cells.Item(89, 14).Value = -37705  # N89: -33953.6675 -> -37705

$ws.Cells.Item(105, 8).Value = 2687.2  # H105: 3234 -> 2687.2
$ws.Cells.Item(105, 9).Value = 2687.2  # I105: 3234 -> 2687.2
$ws.Cells.Item(105, 11).Value = 2687.2  # K105: 3234 -> 2687.2
$ws.Cells.Item(105, 13).Value = -940.1999999999998  # M105: -1487 -> -940.1999999999998

$ws.Cells.Item(107, 8).Value = 1629.1666  # H107: 1495.2858 -> 1629.1666
$ws.Cells.Item(107, 9).Value = 1413.7273  # I107: 1302.6923 -> 1413.7273
$ws.Cells.Item(107, 11).Value = 1413.7273  # K107: 1302.6923 -> 1413.7273
$ws.Cells.Item(107, 13).Value = 506.2727  # M107: 617.3077000000001 -> 506.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3790.8572  # H31: 4707.8 -> 3790.8572
$ws.Cells.Item(31, 9).Value = 5700  # I31: 5949.5 -> 5700
$ws.Cells.Item(31, 10).Value = 3027.2  # J31: 3880 -> 3027.2
$ws.Cells.Item(31, 11).Value = 5700  # K31: 5949.5 -> 5700
$ws.Cells.Item(31, 12).Value = 3027.2  # L31: 3880 -> 3027.2
$ws.Cells.Item(31, 13).Value = -5405  # M31: -5654.5 -> -5405
$ws.Cells.Item(31, 14).Value = -3617.2  # N31: -4470 -> -3617.2

$ws.Cells.Item(34, 8).Value = 3790.8572  # H34: 4707.8 -> 3790.8572
$ws.Cells.Item(34, 9).Value = 5700  # I34: 5949.5 -> 5700
$ws.Cells.Item(34, 10).Value = 3027.2  # J34: 3880 -> 3027.2
$ws.Cells.Item(34, 11).Value = 5700  # K34: 5949.5 -> 5700
$ws.Cells.Item(34, 12).Value = 3027.2  # L34: 3880 -> 3027.2
$ws.Cells.Item(34, 13).Value = -5498  # M34: -5747.5 -> -5498
$ws.Cells.Item(34, 14).Value = -3431.2  # N34: -4284 -> -3431.2

$ws.Cells.Item(68, 8).Value = 0  # H68: 20000 -> 0
$ws.Cells.Item(68, 9).Value = 0  # I68: 20000 -> 0
$ws.Cells.Item(68, 11).Value = 0  # K68: 20000 -> 0
$ws.Cells.Item(68, 13).ClearContents()  # M68: -19251 -> (blank)

$ws.Cells.Item(71, 8).Value = 0  # H71: 20000 -> 0
$ws.Cells.Item(71, 9).Value = 0  # I71: 20000 -> 0
$ws.Cells.Item(71, 11).Value = 0  # K71: 60000 -> 0
$ws.Cells.Item(71, 13).ClearContents()  # M71: -56256 -> (blank)

$ws.Cells.Item(93, 8).Value = 8750  # H93: 13750 -> 8750
$ws.Cells.Item(93, 9).Value = 7500  # I93: 13750 -> 7500
$ws.Cells.Item(93, 10).Value = 10000  # J93: 0 -> 10000
$ws.Cells.Item(93, 11).Value = 7500  # K93: 13750 -> 7500
$ws.Cells.Item(93, 12).Value = 10000  # L93: 0 -> 10000
$ws.Cells.Item(93, 13).Value = -5628  # M93: -11878 -> -5628
$ws.Cells.Item(93, 14).Value = -13744  # N93: (blank) -> -13744

$ws.Cells.Item(96, 8).Value = 15902.857  # H96: 15665 -> 15902.857
$ws.Cells.Item(96, 10).Value = 15902.857  # J96: 15665 -> 15902.857
$ws.Cells.Item(96, 12).Value = 15902.857  # L96: 15665 -> 15902.857
$ws.Cells.Item(96, 14).Value = -21394.857  # N96: -21157 -> -21394.857

$ws.Cells.Item(122, 8).Value = 1106.6875  # H122: 1147.1333 -> 1106.6875
$ws.Cells.Item(122, 9).Value = 737.36365  # I122: 761.1 -> 737.36365
$ws.Cells.Item(122, 11).Value = 2212.09095  # K122: 2283.3 -> 2212.09095
$ws.Cells.Item(122, 13).Value = 237.9090500000002  # M122: 166.6999999999998 -> 237.9090500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 4907.45  # H39: 5254.0586 -> 4907.45
$ws.Cells.Item(39, 10).Value = 4771  # J39: 5113.6875 -> 4771
$ws.Cells.Item(39, 12).Value = 14313  # L39: 15341.0625 -> 14313
$ws.Cells.Item(39, 14).Value = -14901  # N39: -15929.0625 -> -14901

$ws.Cells.Item(44, 8).Value = 0  # H44: 333667.34 -> 0
$ws.Cells.Item(44, 9).Value = 0  # I44: 333667.34 -> 0
$ws.Cells.Item(44, 11).Value = 0  # K44: 1001002.02 -> 0
$ws.Cells.Item(44, 13).ClearContents()  # M44: -1000604.02 -> (blank)

$ws.Cells.Item(107, 8).Value = 100  # H107: 0 -> 100
$ws.Cells.Item(107, 9).Value = 100  # I107: 0 -> 100
$ws.Cells.Item(107, 11).Value = 300  # K107: 0 -> 300
$ws.Cells.Item(107, 13).Value = 1620  # M107: (blank) -> 1620

$ws.Cells.Item(130, 8).Value = 1439.4  # H130: 1299.25 -> 1439.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 11.571428  # H2: 6.4 -> 11.571428
$ws.Cells.Item(2, 9).Value = 8  # I2: 2.4285715 -> 8
$ws.Cells.Item(2, 10).Value = 14.25  # J2: 15.666667 -> 14.25
$ws.Cells.Item(2, 11).Value = 8  # K2: 2.4285715 -> 8
$ws.Cells.Item(2, 12).Value = 14.25  # L2: 15.666667 -> 14.25
$ws.Cells.Item(2, 13).Value = 105  # M2: 110.5714285 -> 105
$ws.Cells.Item(2, 14).Value = -240.25  # N2: -241.666667 -> -240.25

$ws.Cells.Item(102, 8).Value = 669.8333  # H102: 2392.182 -> 669.8333
$ws.Cells.Item(102, 9).Value = 669.8333  # I102: 1835.2222 -> 669.8333
$ws.Cells.Item(102, 10).Value = 0  # J102: 4898.5 -> 0
$ws.Cells.Item(102, 11).Value = 669.8333  # K102: 1835.2222 -> 669.8333
$ws.Cells.Item(102, 12).Value = 0  # L102: 4898.5 -> 0
$ws.Cells.Item(102, 13).ClearContents()  # M102: -213.2221999999999 -> (blank)
$ws.Cells.Item(102, 14).Value = 952.1667  # N102: -8142.5 -> 952.1667

$ws.Cells.Item(104, 8).Value = 100000  # H104: 71500 -> 100000
$ws.Cells.Item(104, 10).Value = 100000  # J104: 71500 -> 100000
$ws.Cells.Item(104, 12).Value = 100000  # L104: 71500 -> 100000
$ws.Cells.Item(104, 14).Value = -106988  # N104: -78488 -> -106988

$ws.Cells.Item(122, 8).Value = 2806.5833  # H122: 2906.5454 -> 2806.5833
$ws.Cells.Item(122, 9).Value = 2806.5833  # I122: 2906.5454 -> 2806.5833
$ws.Cells.Item(122, 11).Value = 8419.749899999999  # K122: 8719.636200000001 -> 8419.749899999999
$ws.Cells.Item(122, 13).Value = -5969.749899999999  # M122: -6269.636200000001 -> -5969.749899999999

$ws.Cells.Item(132, 8).Value = 3692.182  # H132: 3750.35 -> 3692.182
$ws.Cells.Item(132, 9).Value = 3643.6843  # I132: 3706.4119 -> 3643.6843
$ws.Cells.Item(132, 11).Value = 10931.0529  # K132: 11119.2357 -> 10931.0529
$ws.Cells.Item(132, 13).Value = -8401.052899999999  # M132: -8589.235700000001 -> -8401.052899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1150  # H22: 965 -> 1150
$ws.Cells.Item(22, 9).Value = 700  # I22: 656.6667 -> 700
$ws.Cells.Item(22, 10).Value = 1600  # J22: 1273.3334 -> 1600
$ws.Cells.Item(22, 11).Value = 700  # K22: 656.6667 -> 700
$ws.Cells.Item(22, 12).Value = 1600  # L22: 1273.3334 -> 1600
$ws.Cells.Item(22, 13).Value = -405  # M22: -361.6667 -> -405
$ws.Cells.Item(22, 14).Value = -2190  # N22: -1863.3334 -> -2190

$ws.Cells.Item(27, 8).Value = 1150  # H27: 965 -> 1150
$ws.Cells.Item(27, 9).Value = 700  # I27: 656.6667 -> 700
$ws.Cells.Item(27, 10).Value = 1600  # J27: 1273.3334 -> 1600
$ws.Cells.Item(27, 11).Value = 700  # K27: 656.6667 -> 700
$ws.Cells.Item(27, 12).Value = 1600  # L27: 1273.3334 -> 1600
$ws.Cells.Item(27, 13).Value = -593  # M27: -549.6667 -> -593
$ws.Cells.Item(27, 14).Value = -1814  # N27: -1487.3334 -> -1814

$ws.Cells.Item(40, 8).Value = 2997.3635  # H40: 3039 -> 2997.3635
$ws.Cells.Item(40, 9).Value = 2846.5  # I40: 2918.7778 -> 2846.5
$ws.Cells.Item(40, 11).Value = 2846.5  # K40: 2918.7778 -> 2846.5
$ws.Cells.Item(40, 13).Value = -2710.5  # M40: -2782.7778 -> -2710.5

$ws.Cells.Item(48, 8).Value = 30499.8  # H48: 31499.8 -> 30499.8
$ws.Cells.Item(48, 9).Value = 31249.5  # I48: 37499 -> 31249.5
$ws.Cells.Item(48, 11).Value = 31249.5  # K48: 37499 -> 31249.5
$ws.Cells.Item(48, 13).Value = -30588.5  # M48: -36838 -> -30588.5

$ws.Cells.Item(64, 8).Value = 19999.5  # H64: 34999 -> 19999.5
$ws.Cells.Item(64, 10).Value = 19999.5  # J64: 34999 -> 19999.5
$ws.Cells.Item(64, 12).Value = 19999.5  # L64: 34999 -> 19999.5
$ws.Cells.Item(64, 14).Value = -20449.5  # N64: -35449 -> -20449.5

$ws.Cells.Item(67, 8).Value = 19999.5  # H67: 34999 -> 19999.5
$ws.Cells.Item(67, 10).Value = 19999.5  # J67: 34999 -> 19999.5
$ws.Cells.Item(67, 12).Value = 19999.5  # L67: 34999 -> 19999.5
$ws.Cells.Item(67, 14).Value = -21559.5  # N67: -36559 -> -21559.5

$ws.Cells.Item(93, 8).Value = 245.8  # H93: 247.5 -> 245.8
$ws.Cells.Item(93, 9).Value = 245.8  # I93: 247.5 -> 245.8
$ws.Cells.Item(93, 11).Value = 245.8  # K93: 247.5 -> 245.8
$ws.Cells.Item(93, 13).Value = 1002.2  # M93: 1000.5 -> 1002.2

$ws.Cells.Item(106, 8).Value = 0  # H106: 18999 -> 0
$ws.Cells.Item(106, 10).Value = 0  # J106: 18999 -> 0
$ws.Cells.Item(106, 12).ClearContents()  # L106: 18999 -> (blank)
$ws.Cells.Item(106, 14).Value = 0  # N106: -21523 -> 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value = 20000  # H105: 19432.857 -> 20000
$ws.Cells.Item(105, 9).Value = 19000  # I105: 19800 -> 19000
$ws.Cells.Item(105, 10).Value = 21000  # J105: 19371.666 -> 21000
$ws.Cells.Item(105, 11).Value = 19000  # K105: 19800 -> 19000
$ws.Cells.Item(105, 12).Value = 21000  # L105: 19371.666 -> 21000
$ws.Cells.Item(105, 13).Value = -15506  # M105: -16306 -> -15506
$ws.Cells.Item(105, 14).Value = -27988  # N105: -26359.666 -> -27988

$ws.Cells.Item(116, 8).Value = 80000  # H116: 0 -> 80000
$ws.Cells.Item(116, 9).Value = 80000  # I116: 0 -> 80000
$ws.Cells.Item(116, 11).Value = 80000  # K116: 0 -> 80000
$ws.Cells.Item(116, 13).Value = -75411  # M116: (blank) -> -75411

$ws.Cells.Item(131, 8).Value = 29990  # H131: 29998 -> 29990
$ws.Cells.Item(131, 10).Value = 29990  # J131: 29998 -> 29990
$ws.Cells.Item(131, 12).Value = 29990  # L131: 29998 -> 29990
$ws.Cells.Item(131, 14).Value = -40070  # N131: -40078 -> -40070

$ws.Cells.Item(132, 8).Value = 4155.5884  # H132: 4432.25 -> 4155.5884
$ws.Cells.Item(132, 9).Value = 4033.4614  # I132: 4241.3335 -> 4033.4614
$ws.Cells.Item(132, 10).Value = 4552.5  # J132: 5005 -> 4552.5
$ws.Cells.Item(132, 11).Value = 12100.3842  # K132: 12724.0005 -> 12100.3842
$ws.Cells.Item(132, 12).Value = 13657.5  # L132: 15015 -> 13657.5
$ws.Cells.Item(132, 13).Value = -9570.3842  # M132: -10194.0005 -> -9570.3842
$ws.Cells.Item(132, 14).Value = -18717.5  # N132: -20075 -> -18717.5
